$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 values (E4:K4)
$ws.Range("E4").Value = 8322
$ws.Range("F4").Value = 5355
$ws.Range("G4").Value = 5161
$ws.Range("H4").Value = 4888
$ws.Range("I4").Value = 4781
$ws.Range("J4").Value = 5136
$ws.Range("K4").Value = 5425

# Update row 5 values (E5:K5)
$ws.Range("E5").Value = 2562
$ws.Range("F5").Value = 2485
$ws.Range("G5").Value = 2231
$ws.Range("H5").Value = 1987
$ws.Range("I5").Value = 1750
$ws.Range("J5").Value = 2343
$ws.Range("K5").Value = 2710

# Row 5 (E5:K5) style changes from s="15" to s="7" -- copy style from F4 (which is style 7)
$ws.Range("F4:K4").Copy()
$ws.Range("E5:K5").PasteSpecial(-4122)  # xlPasteFormats

# Update selection
$ws.Range("A3").Select()
